# Week 17 data logging + tiebreak-related stat updates
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 17 per-play yardage logs to the long data strings ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "-2 4 4 1 2 5 3 4 4 3 8 9 7 1 2 5 4 7 0 4 0 9 1 7 -1 1 3 9 4 2 7 6 15 1 3 -5 4 -2 -2 4 12 10 4 11 7 8 7 3 4 7 5 14 10 4 2 3 8 4 2 18 2 12 7 5 7 4 -1 2 2 21 5 4 5 4 5 3 12 2 15 24 17 1 1 2 9 3 10 10 1 8 1 2 4 1 12 9 16 3 3 5 4 4 2 10 5 2 3 15 3 6 4 1 9 5 0 -1 23 9 -1 -2 9 5 9 0 7 2 0 6 3 3 17 3 6 1 6 4 0 5 11 4 3 3 3 10 5 2 -1 1 3 0 5 2 2 11 -1 13 18 6 9 9 6 2 5 2 10 -3 4 2 4 4 9 7 5 9 3 8 12 -1 3 4 6 0 -1 1 8 6 -1 7 7 2 6 0 0 6 4 7 7 7 2 4 10 2 -4 -1 -1 2 3 2 5 0 5 1 3 21 1 3 2 0 6 -1 5 3 2 6 -1 4 4 4 3 6 8 33 11 4 3 4 1 5 1 5 0 2 3 3 7 4 14 3 6 4 10 -2 -7 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1 2 1 2 7 5 7 0 3 1 0 0 1 4 9 6 4 2 10 4 2 -2 2 4 2 5 2 5 11 9 7 0 5 0 -1 3 3 1 2 3 0 1 5 0 1 32 0 3 4 0 1 3 8 12 5 5 2 3 2 7 2 6 3 8 3 10 4 6 7 2 7 4 0 1 0 0 1 5 -2 3 1 2 23 -1 6 2 1 8 10 6 8 1 14 3 17 1 19 8 0 5 8 11"
$wsYDS.Range("C2").Value = "2 6 7 -1 6 4 2 4 17 8 18 5 15 6 17 4 2 6 7 6 2 2 -2 5 2 6 5 12 31 9 7 1 9 9 8 0 0 21 3 4 5 20 14 -4 3 4 -2 20 6 2 8 5 9 2 1 11 5 9 -1 3 1 1 2 1 4 -1 9 1 2 3 3 5 -2 -1 5 9 1 -1 11 6 9 6 0 8 1 2 5 11 12 5 12 6 3 3 -4 7 2 2 7 4 2 13 10 0 8 22 2 11 9 9 8 2 -1 3 5 -1 0 0 3 1 4 12 4 6 6 0 8 3 4 8 1 3 5 3 3 2 2 15 1 7 8 6 4 4 11 2 5 3 9 2 3 4 0 8 5 0 3 4 -1 2 4 0 11 1 -1 2 2 4 6 2 1 2 8 5 4 -4 9 3 2 -4 5 2 0 2 -3 11 2 2 6 2 14 5 4 0 3 5 4 10 2 2 4 9 15 2 3 -1 7 7 -2 8 5 3 1 5 5 0 7 8 2 2 4 9 6 4 9 3 5 3 3 1 -1 12 2 -1 9 2 1 3 8 -1 2 5 31 3 1 5 1 5 2 8 3 5 4 4 4 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1 -3 5 1 4 3 4 1 5 3 4 15 3 15 4 3 6 2 4 1 2 7 0 9 4 6 8 4 3 1 8 2 4 8 22 1 1 3 0 2 2 3 2 6 15 2 2 -2 0 2 1 21 5 5 8 6 6 6 12 4 6 3 6 6 5 8 4 2 11 17 2 13 3 6 6 3 1 4 3 3 1 3 4 5 4 1 0 0 0"
$wsYDS.Range("B3").Value = "17 9 11 24 9 8 -6 4 11 6 22 19 -2 9 28 6 7 12 6 11 8 3 14 75 8 4 14 6 7 19 10 33 17 11 14 10 4 14 20 14 5 11 40 5 46 3 8 3 23 13 7 7 15 5 6 6 8 17 -1 5 15 28 7 6 11 12 8 10 2 6 23 18 10 5 5 7 11 8 11 9 36 1 9 2 19 32 4 6 4 9 2 15 7 6 3 14 7 10 5 12 11 44 12 6 15 11 9 7 8 5 8 10 6 26 8 12 8 11 5 4 3 11 8 17 6 8 5 5 4 1 8 6 4 4 11 6 20 9 19 27 49 3 5 7 27 19 17 14 9 10 5 13 9 5 12 13 2 9 0 3 11 5 21 19 6 12 6 24 7 8 6 6 5 10 9 23 8 11 24 12 14 6 4 9 10 4 2 10 3 9 12 7 0 15 5 6 6 7 11 9 17 9 3 11 12 9 18 8 4 9 10 10 15 14 24 2 2 7 10 4 2 1 7 2 7 12 2 25 6 14 10 13 22 9 13 6 7 5 11 22 19 8 8 7 13 12 -1 16 27 3 1 0 9 7 11 7 13 17 2 10 16 32 3 7 3 38 8 30 22 1 10 8 12 0 8 37 13 23 7 11 8 24 14 22 10 12 6 6 4 8 8 8 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4 6 44 23 10 19 19 28 7 9 3 1 7 38 15 5 -3 2 8 16 1 12 14 5 9 20 11 12 -1 17 4 3 6 5 23 10 40 5 3 9 0 6 69 1 5 10 20 7 10 27 14 34 10 11 11 9 13 -1 5 18 18 5 7 3 14 8 50 3 8 9 4 6 9 16 22 8 6 1 29 7 17 24 3 3 53 7 3 2 4 27 7 8 5 11 6 5 6 4 8 1 4"
$wsYDS.Range("C3").Value = "8 14 16 8 44 17 3 9 10 9 43 11 9 30 4 10 3 4 19 11 17 -4 7 20 8 9 7 11 4 4 12 8 27 20 42 19 9 10 6 0 4 3 8 8 11 19 4 20 8 6 16 5 12 7 21 11 5 0 20 43 15 6 16 4 13 22 21 12 10 6 3 5 15 3 7 14 18 6 37 12 1 21 8 7 6 4 14 14 9 18 19 17 9 9 15 23 8 24 35 61 15 53 41 1 -2 11 5 16 16 8 9 13 6 3 7 12 2 7 3 5 6 6 6 4 16 39 10 14 4 -4 -4 9 12 4 -2 7 10 18 3 5 3 24 8 13 22 10 2 20 6 12 9 11 6 5 8 24 5 50 19 3 11 1 12 41 2 5 4 22 4 6 11 6 5 7 16 4 3 3 10 15 5 6 21 35 0 11 -1 14 16 3 8 2 15 5 6 5 1 4 3 27 3 6 7 5 22 37 24 40 11 5 5 7 9 6 12 3 15 -2 8 -3 11 11 2 2 8 11 3 1 5 4 5 16 4 0 6 12 18 14 -2 3 19 5 19 7 11 16 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7 7 8 9 16 14 2 7 19 7 2 5 7 13 10 6 4 6 4 19 19 6 4 9 4 12 12 13 15 4 2 27 7 5 4 17 6 29 8 5 15 15 8 5 14 18 9 5 -3 24 10 0 1 3 5 10 6 2 17 4 8 2 5 6 6 3 3 7 4 6 12 11 4 3 10 72 8 8 9 17 18 6 7 6 4 10 5 69 22 6 4 39 10 5 2 5 9 2 35 30"

# --- OFF sheet: Week 17 offensive totals (Home = row 2, Road = row 3) ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 197
$wsOFF.Range("D2").Value = 14
$wsOFF.Range("F2").Value = 50
$wsOFF.Range("G2").Value = 58
$wsOFF.Range("J2").Value = 27
$wsOFF.Range("O2").Value = 13
$wsOFF.Range("P2").Value = 8
$wsOFF.Range("B3").Value = 13
$wsOFF.Range("C3").Value = 215
$wsOFF.Range("F3").Value = 115
$wsOFF.Range("G3").Value = 60
$wsOFF.Range("H3").Value = 24
$wsOFF.Range("I3").Value = 56
$wsOFF.Range("J3").Value = 59
$wsOFF.Range("L3").Value = 277
$wsOFF.Range("M3").Value = 194
$wsOFF.Range("Q3").Value = 444

# --- DEF sheet: Week 17 defensive totals (Home = row 2, Road = row 3) ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("B2").Value = 5
$wsDEF.Range("C2").Value = 172
$wsDEF.Range("D2").Value = 14
$wsDEF.Range("F2").Value = 66
$wsDEF.Range("G2").Value = 55
$wsDEF.Range("J2").Value = 32
$wsDEF.Range("N2").Value = 11
$wsDEF.Range("O2").Value = 32
$wsDEF.Range("P2").Value = 19
$wsDEF.Range("B3").Value = 12
$wsDEF.Range("C3").Value = 196
$wsDEF.Range("E3").Value = 30
$wsDEF.Range("F3").Value = 104
$wsDEF.Range("G3").Value = 45
$wsDEF.Range("H3").Value = 20
$wsDEF.Range("I3").Value = 58
$wsDEF.Range("J3").Value = 50
$wsDEF.Range("L3").Value = 255
$wsDEF.Range("M3").Value = 176
$wsDEF.Range("Q3").Value = 459

# --- ST sheet: Week 17 special-teams totals + appended per-kick logs ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 90
$wsST.Range("D2").Value = 34
$wsST.Range("F2").Value = 249
$wsST.Range("G2").Value = 233
$wsST.Range("J2").Value = 94
$wsST.Range("K2").Value = 91
$wsST.Range("B3").Value = 59
$wsST.Range("D3").Value = "44 47 54 49 32 51 46 42 34 19 52 41 50 52 43 54 55 44 64 42 49 49 68 31 54 50 37 52 32 41 50 42 35 31"
$wsST.Range("D4").Value = "0 4 0 16 0 0 7 0 0 0 0 0 0 0 0 0 0 0 15 0 0 7 0 0 14 6 0 6 0 0 20 0 3 0"
$wsST.Range("D5").Value = "0 0 0 0 11 8 18 0 0 29 3"
$wsST.Range("B6").Value = "24 20 8 33 3 9 25 30 28 22 32 28 31 33 0 28 23 12 20 14 38 27 39 24 19 20 25 26 30 12"

# --- TURNS sheet: Week 17 turnovers (Road = row 3) ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("D3").Value = 11
$wsTURNS.Range("E3").Value = 11

# --- PEN sheet: Week 17 penalty counts ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("D2").Value = 10
$wsPEN.Range("B3").Value = 22
$wsPEN.Range("D4").Value = 12

Write-Host "Week 17 data logged."
